$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first three mass-point columns (250, 500, 750) so that the
# remaining columns (previously E:K = 1000..2500) shift left to B:H.
$ws.Range("B1:D9").Delete()

# Make room for a brand-new mass-point row (1.8) between the existing
# 1.5 row (row 5) and 2.0 row (row 6 before the insert).
$ws.Range("A6:H6").Insert()

# The inserted row loses the bordered label style on A6; restore it by
# copying the format from the row above (A5), which already carries the
# correct style used for every other mass-point label in column A.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Fill in the new 1.8 row of cross sections.
$ws.Range("A6").Value = "'1.8"
$ws.Range("B6").Value = 0.05132095
$ws.Range("C6").Value = 0.0114448
$ws.Range("D6").Value = 0.00304724
$ws.Range("E6").Value = 0.00091653
$ws.Range("F6").Value = 0.00030121
$ws.Range("G6").Value = 0.00010571
$ws.Range("H6").Value = 0.00003895899
